$p = $ppt.ActivePresentation

# Locate the slide/shape holding the "Singularity developed ... by Berkeley Lab"
# paragraph (Introduction slide, content placeholder) instead of hard-coding
# indices, so the script is resilient to minor deck reordering.
$targetSlide = $null
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $shp = $sl.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -like "*Singularity developed*by Berkeley Lab*") {
                $targetSlide = $sl
                $targetShape = $shp
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Before:
#   "Singularity developed by Berkeley Lab<br>http://singularity.lbl.gov/ "
# After:
#   "Singularity originally developed by Berkeley Lab, now at Sylabs<br>https://sylabs.io/singularity/"

# 1) "Singularity developed " -> "Singularity originally " (keeps "developed " as
#    its own trailing run, same as the source split).
$tr.Characters(1, 12).Text = "Singularity originally "

# 2) Expand the trailing "Lab" of "by Berkeley Lab" into "Lab, now at Sylabs".
$tr.Characters(46, 3).Text = "Lab, now at Sylabs"

# 3) Split that into "Lab, now at " + "Sylabs" runs.
$tr.Characters(46, 12).Text = "Lab, now at "

# 4) Replace the old hyperlink display text (and the stray trailing "/" and
#    space runs that followed it) with the new Sylabs URL, in a single run.
$tr.Characters(65, 28).Text = "https://sylabs.io/singularity/"

# 5) Point the existing hyperlink (rId2) at the new address so the relationship
#    target actually matches the new display text.
$urlRange = $tr.Characters(65, 30)
$urlRange.ActionSettings(1).Hyperlink.Address = "https://sylabs.io/singularity/"
